# Update the division problems in the single table of the worksheet.
# Cells are addressed directly by (row, column) rather than via a global
# Find/Replace, because some old values equal other new values
# (e.g. "479÷7=" is both the new text of Cell(1,1) and the old text of
# Cell(9,5)), which would make a naive text-based replace ambiguous.
$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Row 1 (problems 1-5)
$tbl.Cell(1, 1).Range.Text = "479÷7="
$tbl.Cell(1, 2).Range.Text = "220÷8="
$tbl.Cell(1, 3).Range.Text = "589÷5="
$tbl.Cell(1, 4).Range.Text = "705÷3="
$tbl.Cell(1, 5).Range.Text = "458÷8="

# Row 5 (problems 6-10)
$tbl.Cell(5, 1).Range.Text = "394÷6="
$tbl.Cell(5, 2).Range.Text = "725÷2="
$tbl.Cell(5, 3).Range.Text = "923÷9="
$tbl.Cell(5, 4).Range.Text = "339÷3="
$tbl.Cell(5, 5).Range.Text = "816÷7="

# Row 9 (problems 11-15)
$tbl.Cell(9, 1).Range.Text = "309÷6="
$tbl.Cell(9, 2).Range.Text = "345÷9="
$tbl.Cell(9, 3).Range.Text = "676÷9="
$tbl.Cell(9, 4).Range.Text = "849÷5="
$tbl.Cell(9, 5).Range.Text = "958÷3="

# Row 13 (problems 16-20)
$tbl.Cell(13, 1).Range.Text = "660÷3="
$tbl.Cell(13, 2).Range.Text = "823÷3="
$tbl.Cell(13, 3).Range.Text = "813÷5="
$tbl.Cell(13, 4).Range.Text = "825÷5="
$tbl.Cell(13, 5).Range.Text = "281÷2="

# Row 17 (problems 21-25)
$tbl.Cell(17, 1).Range.Text = "387÷6="
$tbl.Cell(17, 2).Range.Text = "500÷2="
$tbl.Cell(17, 3).Range.Text = "771÷9="
$tbl.Cell(17, 4).Range.Text = "980÷5="
$tbl.Cell(17, 5).Range.Text = "183÷6="
